$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 182, pushing existing rows 182-269 down to 184-271
$ws.Rows.Item(182).Insert()
$ws.Rows.Item(182).Insert()

# Populate new row 182 (Primera) - new weekly entry
$ws.Range("A182").Value = 11
$ws.Range("B182").Value = "Vega Monumental Concepción"
$ws.Range("C182").Value = "Bíobío"
$ws.Range("D182").Value = 44489
$ws.Range("E182").Value = 8
$ws.Range("F182").Value = 100112020
$ws.Range("G182").Value = "Tomate"
$ws.Range("H182").Value = "Larga vida"
$ws.Range("I182").Value = "Primera"
$ws.Range("J182").Value = 2000
$ws.Range("K182").Value = 6000
$ws.Range("L182").Value = 6500
$ws.Range("M182").Value = 6250
$ws.Range("N182").Value = "$/caja 10 kilos"
$ws.Range("O182").Value = "Región de Arica y Parinacota"
$ws.Range("P182").Value = 625
$ws.Range("Q182").Value = 10
$ws.Range("R182").Value = "Hortaliza"

# Populate new row 183 (Segunda) - new weekly entry
$ws.Range("A183").Value = 11
$ws.Range("B183").Value = "Vega Monumental Concepción"
$ws.Range("C183").Value = "Bíobío"
$ws.Range("D183").Value = 44489
$ws.Range("E183").Value = 8
$ws.Range("F183").Value = 100112020
$ws.Range("G183").Value = "Tomate"
$ws.Range("H183").Value = "Larga vida"
$ws.Range("I183").Value = "Segunda"
$ws.Range("J183").Value = 1000
$ws.Range("K183").Value = 5000
$ws.Range("L183").Value = 5000
$ws.Range("M183").Value = 5000
$ws.Range("N183").Value = "$/caja 10 kilos"
$ws.Range("O183").Value = "Región de Arica y Parinacota"
$ws.Range("P183").Value = 500
$ws.Range("Q183").Value = 10
$ws.Range("R183").Value = "Hortaliza"
